# "Fruta / hortaliza, semanal" — weekly data refresh.
# A new weekly price observation is inserted at row 141 (Región de Arica y
# Parinacota, $/caja 60 unidades, Pepino ensalada), pushing every existing
# row from 141 down through 217 to 142 through 218. The sheet's dimension
# grows from A1:R217 to A1:R218.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 141, shifting rows 141:217 down
# to 142:218 (this is what the diff shows: row N's old content ends up in
# row N+1, all the way down, with a brand new row appearing at 141 and the
# sheet growing by one row at the bottom).
$ws.Rows(141).Insert()

# Populate the newly inserted row 141 with the new observation.
$ws.Range("A141").Value = 5
$ws.Range("B141").Value = "Macroferia Regional de Talca"
$ws.Range("C141").Value = "Maule"
$ws.Range("D141").Value = 44460
$ws.Range("E141").Value = 7
$ws.Range("F141").Value = 100112043
$ws.Range("G141").Value = "Pepino ensalada"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 400
$ws.Range("K141").Value = 16000
$ws.Range("L141").Value = 16000
$ws.Range("M141").Value = 16000
$ws.Range("N141").Value = "$/caja 60 unidades"
$ws.Range("O141").Value = "Región de Arica y Parinacota"
$ws.Range("P141").Value = 267
$ws.Range("Q141").Value = 60
$ws.Range("R141").Value = "Hortaliza"
